$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in C1: "Cólùm 3" -> "Cólùmn 3" (B2 "Test string" stays the same)
$ws.Range("C1").Value = "Cólùmn 3"

# Move the active selection from A2 to C1
$ws.Range("C1").Select()
